$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the old <w:bookmarkStart/bookmarkEnd w:name="_GoBack"/> pair that
#    sat inside the numbered-list paragraph ("Create a new repositories...").
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete() | Out-Null
Write-Host "Removed old _GoBack bookmark"

# ---------------------------------------------------------------------------
# 2) Merge the three runs that together spell out
#    "https://github.com/monacog/" + "git-course" + ".git" into a single run
#    (first occurrence only - the git-clone help text block).
# ---------------------------------------------------------------------------
$urlRng = $d.Content
$urlRng.Find.Execute("https://github.com/monacog/", $true, $false, $false, $false, $false, `
                      $true, 1, $false, $null, 0) | Out-Null
$urlRng.InsertAfter("git-course.git") | Out-Null

$tailRng = $d.Content
$tailRng.Start = $urlRng.End
$tailRng.Find.Execute("git-course.git", $true, $false, $false, $false, $false, `
                       $true, 1, $false, $null, 0) | Out-Null
$tailRng.Delete() | Out-Null
Write-Host "Merged git-course clone URL runs"

# ---------------------------------------------------------------------------
# 3) Add the new "Formato de escritura : MARKDOWN" paragraph run, and move the
#    _GoBack bookmark to sit right after it (end of document).
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastRng = $lastPara.Range
$newParaXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' " + `
    "xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml' " + `
    "w14:paraId='348CC9CD' w14:textId='18269691' w:rsidR='009458B3' w:rsidRDefault='009458B3'>" + `
    "<w:pPr><w:rPr><w:lang w:val='es-MX'/></w:rPr></w:pPr>" + `
    "<w:r><w:rPr><w:lang w:val='es-MX'/></w:rPr><w:t>Formato de escritura : MARKDOWN</w:t></w:r>" + `
    "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" + `
    "</w:p>"
$lastRng.InsertXML($newParaXml) | Out-Null
Write-Host "Added Markdown format note paragraph with relocated _GoBack bookmark"
